# edit.ps1 - applies the changes described by the diff to resourcesData.xlsx
#
# Summary of the target edit:
#  - basicGeom (sheet1): B5:B13 values reset from 1..9 to 0
#  - basicGeom (sheet1): sheet view no longer frozen/scrolled to A4; becomes
#    the active (selected) tab with the selection moved to E17
#  - basicGeom (sheet1): page setup changed to paper size 9 (A4),
#    portrait orientation
#  - procMap (sheet2): no longer the selected tab; selection moves to A8

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # basicGeom
$ws2 = $wb.Worksheets.Item(2)   # procMap

# --- basicGeom: zero out the B5:B13 column values -----------------------
$ws1.Range("B5:B13").Value = 0

# --- basicGeom: page setup (A4 paper, portrait) --------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- selection / active-tab changes --------------------------------------
# Select procMap's new cell first, then basicGeom last so basicGeom ends up
# as the active sheet/tab (matching tabSelected moving to basicGeom).
$ws2.Range("A8").Select() | Out-Null
$ws1.Range("E17").Select() | Out-Null
